# Auto-generated Excel COM-interop script applying the Anima_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 700.4
$ws.Range("I98").Value = 292.16666
$ws.Range("J98").Value = 2333.3333
$ws.Range("K98").Value = 292.16666
$ws.Range("L98").Value = 2333.3333
$ws.Range("M98").Value = 1205.83334
$ws.Range("N98").Value = -5329.3333

$ws.Range("H122").Value = 700.4
$ws.Range("I122").Value = 292.16666
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 876.4999799999999
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = 1573.50002
$ws.Range("N122").Value = -11899.9999

$ws.Range("H132").Value = 4887.2256
$ws.Range("I132").Value = 4926.304
$ws.Range("J132").Value = 4774.875
$ws.Range("K132").Value = 14778.912
$ws.Range("L132").Value = 14324.625
$ws.Range("M132").Value = -12248.912
$ws.Range("N132").Value = -19384.625

$ws.Range("H137").Value = 2980.1777
$ws.Range("I137").Value = 2118.2646
$ws.Range("J137").Value = 5644.273
$ws.Range("K137").Value = 6354.793799999999
$ws.Range("L137").Value = 16932.819
$ws.Range("M137").Value = -3804.793799999999
$ws.Range("N137").Value = -22032.819

$ws.Range("H138").Value = 2044.86
$ws.Range("J138").Value = 2269.2588
$ws.Range("L138").Value = 6807.776400000001
$ws.Range("N138").Value = -17087.7764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2428.2258
$ws.Range("I61").Value = 1608.75
$ws.Range("J61").Value = 3918.182
$ws.Range("K61").Value = 1608.75
$ws.Range("L61").Value = 3918.182
$ws.Range("M61").Value = -1396.75
$ws.Range("N61").Value = -4342.182

$ws.Range("H74").Value = 1983.2245
$ws.Range("I74").Value = 1560.5
$ws.Range("J74").Value = 2461.087
$ws.Range("K74").Value = 1560.5
$ws.Range("L74").Value = 2461.087
$ws.Range("M74").Value = -686.5
$ws.Range("N74").Value = -4209.087

$ws.Range("H77").Value = 1983.2245
$ws.Range("I77").Value = 1560.5
$ws.Range("J77").Value = 2461.087
$ws.Range("K77").Value = 7802.5
$ws.Range("L77").Value = 12305.435
$ws.Range("M77").Value = -3434.5
$ws.Range("N77").Value = -21041.435

$ws.Range("H110").Value = 1076.6
$ws.Range("I110").Value = 1082.0714
$ws.Range("K110").Value = 1082.0714
$ws.Range("M110").Value = 962.9286

$ws.Range("H132").Value = 3626.9778
$ws.Range("I132").Value = 3546
$ws.Range("J132").Value = 3728.2
$ws.Range("K132").Value = 10638
$ws.Range("L132").Value = 11184.6
$ws.Range("M132").Value = -8108
$ws.Range("N132").Value = -16244.6

$ws.Range("H136").Value = 2428.2258
$ws.Range("I136").Value = 1608.75
$ws.Range("J136").Value = 3918.182
$ws.Range("K136").Value = 4826.25
$ws.Range("L136").Value = 11754.546
$ws.Range("M136").Value = -2276.25
$ws.Range("N136").Value = -16854.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 64750
$ws.Range("J21").Value = 64750
$ws.Range("L21").Value = 64750
$ws.Range("N21").Value = -65222

$ws.Range("H96").Value = 18204.666
$ws.Range("I96").Value = 12107
$ws.Range("J96").Value = 30400
$ws.Range("K96").Value = 12107
$ws.Range("L96").Value = 30400
$ws.Range("M96").Value = -9361
$ws.Range("N96").Value = -35892

$ws.Range("H97").Value = 31330
$ws.Range("J97").Value = 36995
$ws.Range("L97").Value = 36995
$ws.Range("N97").Value = -38977

$ws.Range("H106").Value = 68132.2
$ws.Range("J106").Value = 68132.2
$ws.Range("L106").Value = 68132.2
$ws.Range("N106").Value = -70656.2

$ws.Range("H134").Value = 2133.311
$ws.Range("I134").Value = 1811.2424
$ws.Range("J134").Value = 3019
$ws.Range("K134").Value = 5433.7272
$ws.Range("L134").Value = 9057
$ws.Range("M134").Value = -2898.7272
$ws.Range("N134").Value = -14127

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11745.381
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 11745.381
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 11745.381
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -12335.381

$ws.Range("H34").Value = 11745.381
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 11745.381
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 11745.381
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -12149.381

$ws.Range("H58").Value = 1396.2693
$ws.Range("I58").Value = 1094.9231
$ws.Range("J58").Value = 1697.6154
$ws.Range("K58").Value = 1094.9231
$ws.Range("L58").Value = 1697.6154
$ws.Range("M58").Value = -891.9231
$ws.Range("N58").Value = -2103.6154

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H132").Value = 8335249.5
$ws.Range("I132").Value = 1394.6154
$ws.Range("J132").Value = 23812410
$ws.Range("K132").Value = 4183.8462
$ws.Range("L132").Value = 71437230
$ws.Range("M132").Value = -1653.8462
$ws.Range("N132").Value = -71442290

$ws.Range("H134").Value = 2305.6206
$ws.Range("I134").Value = 1503.7142
$ws.Range("J134").Value = 4410.625
$ws.Range("K134").Value = 4511.142599999999
$ws.Range("L134").Value = 13231.875
$ws.Range("M134").Value = -1976.142599999999
$ws.Range("N134").Value = -18301.875

$ws.Range("H136").Value = 1396.2693
$ws.Range("I136").Value = 1094.9231
$ws.Range("J136").Value = 1697.6154
$ws.Range("K136").Value = 3284.7693
$ws.Range("L136").Value = 5092.8462
$ws.Range("M136").Value = -734.7692999999999
$ws.Range("N136").Value = -10192.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2642.7646
$ws.Range("I132").Value = 2570.8
$ws.Range("K132").Value = 23137.2
$ws.Range("M132").Value = -20607.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5711.5884
$ws.Range("I70").Value = 5691.8
$ws.Range("K70").Value = 5691.8
$ws.Range("M70").Value = -5421.8

$ws.Range("H73").Value = 5711.5884
$ws.Range("I73").Value = 5691.8
$ws.Range("K73").Value = 5691.8
$ws.Range("M73").Value = -4755.8

$ws.Range("H86").Value = 49825.715
$ws.Range("J86").Value = 49825.715
$ws.Range("L86").Value = 49825.715
$ws.Range("N86").Value = -52197.715

$ws.Range("H89").Value = 49825.715
$ws.Range("J89").Value = 49825.715
$ws.Range("L89").Value = 149477.145
$ws.Range("N89").Value = -161333.145

$ws.Range("H132").Value = 3044.04
$ws.Range("I132").Value = 2733.1428
$ws.Range("J132").Value = 3439.7273
$ws.Range("K132").Value = 8199.428400000001
$ws.Range("L132").Value = 10319.1819
$ws.Range("M132").Value = -5669.428400000001
$ws.Range("N132").Value = -15379.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 865.9167
$ws.Range("I16").Value = 865.9167
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 865.9167
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -695.9167
$ws.Range("N16").ClearContents()

$ws.Range("H22").Value = 7954.7036
$ws.Range("I22").Value = 707.3333
$ws.Range("J22").Value = 10025.381
$ws.Range("K22").Value = 707.3333
$ws.Range("L22").Value = 10025.381
$ws.Range("M22").Value = -412.3333
$ws.Range("N22").Value = -10615.381

$ws.Range("H27").Value = 7954.7036
$ws.Range("I27").Value = 707.3333
$ws.Range("J27").Value = 10025.381
$ws.Range("K27").Value = 707.3333
$ws.Range("L27").Value = 10025.381
$ws.Range("M27").Value = -600.3333
$ws.Range("N27").Value = -10239.381

$ws.Range("H46").Value = 1137.5
$ws.Range("J46").Value = 1166.6666
$ws.Range("L46").Value = 1166.6666
$ws.Range("N46").Value = -1542.6666

$ws.Range("H132").Value = 3987.1738
$ws.Range("I132").Value = 3407.923
$ws.Range("K132").Value = 10223.769
$ws.Range("M132").Value = -7693.769

$ws.Range("H136").Value = 7938803
$ws.Range("I136").Value = 2140
$ws.Range("K136").Value = 6420
$ws.Range("M136").Value = -3870

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 32000
$ws.Range("J56").Value = 24000
$ws.Range("L56").Value = 24000
$ws.Range("N56").Value = -25428

$ws.Range("H69").Value = 17085.5
$ws.Range("J69").Value = 17085.5
$ws.Range("L69").Value = 17085.5
$ws.Range("N69").Value = -18583.5

$ws.Range("H72").Value = 17085.5
$ws.Range("J72").Value = 17085.5
$ws.Range("L72").Value = 51256.5
$ws.Range("N72").Value = -58744.5

$ws.Range("H82").Value = 26142.562
$ws.Range("J82").Value = 26142.562
$ws.Range("L82").Value = 26142.562
$ws.Range("N82").Value = -26908.562

$ws.Range("H85").Value = 26142.562
$ws.Range("J85").Value = 26142.562
$ws.Range("L85").Value = 26142.562
$ws.Range("N85").Value = -28794.562

$ws.Range("H132").Value = 3088492.5
$ws.Range("I132").Value = 2354.3635
$ws.Range("K132").Value = 7063.0905
$ws.Range("M132").Value = -4533.0905

$ws.Range("H136").Value = 4248.722
$ws.Range("I136").Value = 4289.3335
$ws.Range("K136").Value = 12868.0005
$ws.Range("M136").Value = -10318.0005
